# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit diff:
#  - Updates the "last updated" timestamp string
#  - Refreshes numeric stats for several countries
#  - Three country pairs swapped ranking order (new data made one country's
#    totals overtake its neighbour's), so the country labels at those two
#    adjacent rows must be exchanged while their numeric rows receive the
#    appropriate (new vs. previously-existing) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the "Datos actualizados..." timestamp banner (A1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 26 de Agosto de 2020 a las 14:14"

# ---------------------------------------------------------------------
# 2. Simple numeric refreshes (country identity/row unchanged)
# ---------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value2 = 5956360
$ws.Range("C4").Value2 = 632
$ws.Range("E4").Value2 = 2519191
$ws.Range("G4").Value2 = 26
$ws.Range("H4").Value2 = 182430

# Row 39 - Oman
$ws.Range("B39").Value2 = 84818
$ws.Range("C39").Value2 = 166
$ws.Range("D39").Value2 = 79409
$ws.Range("E39").Value2 = 4763
$ws.Range("G39").Value2 = 4
$ws.Range("H39").Value2 = 646

# Row 75 - Estado de Palestina
$ws.Range("B75").Value2 = 20155
$ws.Range("C75").Value2 = 477
$ws.Range("D75").Value2 = 13929
$ws.Range("E75").Value2 = 6091
$ws.Range("G75").Value2 = 2
$ws.Range("H75").Value2 = 135

# Row 77 - Bosnia y Herzegovina
$ws.Range("B77").Value2 = 18609
$ws.Range("C77").Value2 = 283
$ws.Range("D77").Value2 = 12336
$ws.Range("E77").Value2 = 5702
$ws.Range("G77").Value2 = 11
$ws.Range("H77").Value2 = 571

# Row 80 - Dinamarca
$ws.Range("B80").Value2 = 16537
$ws.Range("C80").Value2 = 57
$ws.Range("D80").Value2 = 14603
$ws.Range("E80").Value2 = 1311

# Row 82 - Madagascar
$ws.Range("B82").Value2 = 14554
$ws.Range("C82").Value2 = 79
$ws.Range("D82").Value2 = 13582
$ws.Range("E82").Value2 = 791
$ws.Range("G82").Value2 = 3
$ws.Range("H82").Value2 = 181

# Row 139 - Islandia
$ws.Range("B139").Value2 = 2082
$ws.Range("C139").Value2 = 5
$ws.Range("D139").Value2 = 1957
$ws.Range("E139").Value2 = 115

# Row 153 - Burkina Faso
$ws.Range("B153").Value2 = 1352
$ws.Range("C153").Value2 = 14
$ws.Range("D153").Value2 = 1058
$ws.Range("E153").Value2 = 239

# Row 162 - Vietnam
$ws.Range("B162").Value2 = 1034
$ws.Range("C162").Value2 = 5
$ws.Range("D162").Value2 = 632
$ws.Range("E162").Value2 = 373
$ws.Range("G162").Value2 = 2
$ws.Range("H162").Value2 = 29

# ---------------------------------------------------------------------
# 3. Country-pair ranking swaps.
#
# New data pushed each of these countries above its former neighbour, so
# the two adjacent rows swap which country they display. Writing the two
# swapped labels directly (A41="Kuwait" then A42="Rumania") would make
# Excel temporarily see "Rumania" with zero references and reroute/merge
# the shared-string table unexpectedly, so the first cell is parked on a
# unique placeholder value while the second cell's label is settled, and
# only then is the first cell set to its final label.
# ---------------------------------------------------------------------

# Rows 41/42: Rumania/Kuwait -> Kuwait/Rumania
$ws.Range("A41").Value2 = "__SWAP_PLACEHOLDER_1__"
$ws.Range("A42").Value2 = "Rumania"
$ws.Range("A41").Value2 = "Kuwait"
# Row 41 (now Kuwait) gets the new Kuwait figures
$ws.Range("B41").Value2 = 82271
$ws.Range("C41").Value2 = 698
$ws.Range("D41").Value2 = 73906
$ws.Range("E41").Value2 = 7844
$ws.Range("G41").Value2 = 2
$ws.Range("H41").Value2 = 521
# Row 42 (now Rumania) carries over the previous Rumania figures
$ws.Range("B42").Value2 = 81646
$ws.Range("C42").Value2 = 1256
$ws.Range("D42").Value2 = 36286
$ws.Range("E42").Value2 = 41939
$ws.Range("G42").Value2 = 54
$ws.Range("H42").Value2 = 3421

# Rows 66/67: Moldavia/Nepal -> Nepal/Moldavia
$ws.Range("A66").Value2 = "__SWAP_PLACEHOLDER_2__"
$ws.Range("A67").Value2 = "Moldavia"
$ws.Range("A66").Value2 = "Nepal"
# Row 66 (now Nepal) gets the new Nepal figures
$ws.Range("B66").Value2 = 34418
$ws.Range("C66").Value2 = 885
$ws.Range("D66").Value2 = 19504
$ws.Range("E66").Value2 = 14739
$ws.Range("G66").Value2 = 11
$ws.Range("H66").Value2 = 175
# Row 67 (now Moldavia) carries over the previous Moldavia figures
$ws.Range("B67").Value2 = 34358
$ws.Range("D67").Value2 = 23869
$ws.Range("E67").Value2 = 9529
$ws.Range("H67").Value2 = 960

# Rows 96/97: Albania/Croacia -> Croacia/Albania
$ws.Range("A96").Value2 = "__SWAP_PLACEHOLDER_3__"
$ws.Range("A97").Value2 = "Albania"
$ws.Range("A96").Value2 = "Croacia"
# Row 96 (now Croacia) gets the new Croacia figures
$ws.Range("B96").Value2 = 8888
$ws.Range("C96").Value2 = 358
$ws.Range("D96").Value2 = 6362
$ws.Range("E96").Value2 = 2351
$ws.Range("H96").Value2 = 175
# Row 97 (now Albania) carries over the previous Albania figures
$ws.Range("B97").Value2 = 8759
$ws.Range("D97").Value2 = 4530
$ws.Range("E97").Value2 = 3970
$ws.Range("H97").Value2 = 259
